$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the "Transaction Number" column (D) for the table rows that were
# entered/imported from the automation app.
$ws.Range("D2").Value = 247181
$ws.Range("D3").Value = 247182
$ws.Range("D4").Value = 247183
$ws.Range("D5").Value = 247184
$ws.Range("D6").Value = 247185

# Re-fit the data columns so their width matches the (now wider) content,
# mirroring Excel's own "best fit" column sizing for the table.
$ws.Columns("A:C").AutoFit()

# Leave the selection where the user last clicked before saving.
$ws.Range("C17").Select() | Out-Null
